# Release MHD 4.2.2 close #419
#
# Update the "Metadata" worksheet of the MHD List Types ValueSet workbook:
#   - bump the Version from 4.2.1 to 4.2.2
#   - bump the Date to the new release timestamp
#   - replace the placeholder "No display for ContactDetail" contact rows
#     with the three real contact detail lines now rendered for the
#     publisher's ContactDetail (telecom URL, telecom email, and the
#     named contact with its email)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "4.2.2"
$ws.Range("B8").Value = "2024-05-18T12:39:23-05:00"
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"
